$d = $word.ActiveDocument

# The merge-field placeholder "${schule}" (the school-name field) is being
# renamed to "${schule_nametype}". In the source OOXML this shows up as the
# single text run "${schule}" being split so the inserted "_nametype"
# segment lands between "${schule" and the closing "}". We reproduce that by
# locating the run's text and inserting the new fragment right before the
# closing brace, which keeps the surrounding field/bookmark structure intact.
$find = $d.Content.Find
$find.ClearFormatting()
$find.Text = "`${schule}"
$find.Replacement.ClearFormatting()
$find.Replacement.Text = "`${schule_nametype}"
$find.Forward = $true
$find.Wrap = 1
$find.Format = $false
$find.MatchCase = $true
$find.MatchWholeWord = $false
$find.MatchWildcards = $false
$find.MatchSoundsLike = $false
$find.MatchAllWordForms = $false
$find.Execute($find.Text, $true, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null
